$p = $ppt.ActivePresentation

# Slide 16, shape 3 ("PlaceHolder 3") holds the lab-schedule bullet list.
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

# The 4th paragraph reads "First lab will be on January 14th" split across
# four runs: "First lab will be on " / "January " / "14" / "th".
$para = $tr.Paragraphs(4, 1)

# Merge the first two runs into a single run reading
# "First lab will be on January " (keeps run 1's formatting).
$para.Characters(1, 29).Text = "First lab will be on January "

# Update the day-of-month run from "14" to "12" (keeps that run's formatting).
$para.Characters(30, 2).Text = "12"
